$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 10.89835490497452
$ws.Range("A3").Value = 24.62750397673582
$ws.Range("A4").Value = 20.88289491432329
$ws.Range("A5").Value = 22.1714314217609
$ws.Range("A6").Value = 22.2459716575656
$ws.Range("A7").Value = 19.8627667704732
$ws.Range("A8").Value = 17.51806114281132
$ws.Range("A9").Value = 19.88765161502965
$ws.Range("A10").Value = 13.87313780430065
$ws.Range("A11").Value = 20.1905144245213
$ws.Range("A12").Value = 22.20159577114279
$ws.Range("A13").Value = 14.61544706729242
$ws.Range("A14").Value = 14.2277478471413
$ws.Range("A15").Value = 16.45236482271528
$ws.Range("A16").Value = 14.37093776810019
$ws.Range("A17").Value = 13.49500675559565
$ws.Range("A18").Value = 18.63777464806455
$ws.Range("A19").Value = 10.57455121237427
$ws.Range("A20").Value = 14.51716305219858
$ws.Range("A21").Value = 12.50492244635581
$ws.Range("A22").Value = 11.99880493412905
$ws.Range("A23").Value = 13.74415421457692
$ws.Range("A24").Value = 11.30057257976017
$ws.Range("A25").Value = 11.73164096260487
$ws.Range("A26").Value = 8.036845902372278
$ws.Range("A27").Value = 9.994628589008698
$ws.Range("A28").Value = 13.53282529787364
$ws.Range("A29").Value = 8.377997630179237
$ws.Range("A30").Value = 8.432203183078542
$ws.Range("A31").Value = 5.461042980876613
$ws.Range("A32").Value = 8.86906268077874
$ws.Range("A33").Value = 9.915064404167424
$ws.Range("A34").Value = 10.39602741959368
$ws.Range("A35").Value = 11.28706645135679
$ws.Range("A36").Value = 8.15975534162385
$ws.Range("A37").Value = 7.973164742951411
$ws.Range("A38").Value = 7.086360762773353
$ws.Range("A39").Value = 7.841572175589988
$ws.Range("A40").Value = 5.981930266154478
$ws.Range("A41").Value = 5.791030265471221
$ws.Range("A42").Value = 6.26027017159663
$ws.Range("A43").Value = 9.456097019158648
$ws.Range("A44").Value = 8.109963575079547
$ws.Range("A45").Value = 10.88650529577026
$ws.Range("A46").Value = 12.40542978044638
$ws.Range("A47").Value = 8.606416124972782
$ws.Range("A48").Value = 8.9217683513466
$ws.Range("A49").Value = 6.944412901889081
$ws.Range("A50").Value = 8.950332879895939
$ws.Range("A51").Value = 6.932683172820759
$ws.Range("A52").Value = 8.346943068934337
$ws.Range("A53").Value = 9.373061270529092
$ws.Range("A54").Value = 6.675810945434108
$ws.Range("A55").Value = 7.201672055830755
$ws.Range("A56").Value = 7.462466999753644
$ws.Range("A57").Value = 9.534936574661771
$ws.Range("A58").Value = 7.957001635313219
$ws.Range("A59").Value = 7.69888144684748
$ws.Range("A60").Value = 8.124326528385666
$ws.Range("A61").Value = 7.219530104558743
$ws.Range("A62").Value = 6.566454011678218
$ws.Range("A63").Value = 4.391146095249752
$ws.Range("A64").Value = 3.93816971401219
$ws.Range("A65").Value = 8.456032965627941
$ws.Range("A66").Value = 4.234850931000324
$ws.Range("A67").Value = 7.087035791471692
$ws.Range("A68").Value = 4.752798452809401
$ws.Range("A69").Value = 4.227564117610541
$ws.Range("A70").Value = 6.615625445471693
$ws.Range("A71").Value = 8.602491958863084
$ws.Range("A72").Value = 5.223301369830125
$ws.Range("A73").Value = 7.229241319596071
$ws.Range("A74").Value = 4.162270926113592
$ws.Range("A75").Value = 9.234831465217525
$ws.Range("A76").Value = 8.061591374610117
$ws.Range("A77").Value = 7.853951870358628
$ws.Range("A78").Value = 6.939010895201619
$ws.Range("A79").Value = 7.670603352320427
$ws.Range("A80").Value = 9.455493754549593
$ws.Range("A81").Value = 7.257951807573249
$ws.Range("A82").Value = 9.421277742469471
$ws.Range("A83").Value = 7.69714408247188
$ws.Range("A84").Value = 7.13286838164791
$ws.Range("A85").Value = 9.765923334843393
$ws.Range("A86").Value = 7.165612879813466
